# Update the ticker-list table on Sheet1.
# Column A (row index 0..11) is unchanged; only columns B..F of rows 2-13
# get new values (some now blank), and rows 14-18 are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..F, rows 2..13. Empty string clears the cell.
$values = @{
    2  = @("NSE:ATUL",       "NSE:A2ZINFRA",  "NSE:BANDHANBNK", "",  "NSE:MARICO")
    3  = @("NSE:BFSI",       "NSE:ABFRL",     "NSE:IRFC",       "",  "")
    4  = @("NSE:CESC",       "NSE:DOLATALGO", "NSE:RVNL",       "",  "")
    5  = @("NSE:CUB",        "NSE:GNFC",      "",               "",  "")
    6  = @("NSE:EQUITASBNK", "NSE:IFBAGRO",   "",               "",  "")
    7  = @("NSE:FINEORG",    "NSE:IOLCP",     "",               "",  "")
    8  = @("NSE:HDFCLIQUID", "NSE:KALYANKJIL","",               "",  "")
    9  = @("NSE:INGERRAND",  "NSE:KECL",      "",               "",  "")
    10 = @("NSE:MARICO",     "NSE:MAFANG",    "",               "",  "")
    11 = @("NSE:MODISONLTD", "NSE:MOIL",      "",               "",  "")
    12 = @("",               "NSE:NITCO",     "",               "",  "")
    13 = @("",               "NSE:PIIND",     "",               "",  "")
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals[0]
    $ws.Range("C$row").Value = $rowVals[1]
    $ws.Range("D$row").Value = $rowVals[2]
    $ws.Range("E$row").Value = $rowVals[3]
    $ws.Range("F$row").Value = $rowVals[4]
}

# Remove rows 14-18 (the table shrinks from 16 data rows to 11).
$ws.Range("A14:F18").EntireRow.Delete()
